$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = 'orig_sharpe'
$ws.Range("K1").Value = 'orig_n_trades'
$ws.Range("L1").Value = 'orig_commission'
$ws.Range("M1").Value = 'orig_initial_capital'
$ws.Range("N1").Value = 'orig_final_capital'
$ws.Range("O1").Value = 'improved_sharpe'
$ws.Range("P1").Value = 'improved_n_trades'
$ws.Range("Q1").Value = 'improved_commission'
$ws.Range("R1").Value = 'improved_initial_capital'
$ws.Range("S1").Value = 'improved_final_capital'
$ws.Range("T1").Value = 'orig_pct_diff_entry_threshold_pct'
$ws.Range("U1").Value = 'orig_pct_diff_exit_threshold_pct'

# Match header style (bold, bordered, centered) used by A1:I1
$ws.Range("A1").Copy()
$ws.Range("J1:U1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 2 (orig = Bollinger band variant) ----
$ws.Range("A2").Value = 45787.60035719907
$b2 = @'
content="### 1️⃣ Explanation: 
Mean Reversion Strategy for TSLA: Buy when 5-min Close price falls below Lower Bollinger Band (2 std devs), Sell when above Upper Bollinger Band, otherwise Hold.
### 2️⃣ `add_signal` Function

'@
$ws.Range("B2").Value = $b2
$c2 = @'

def add_signal(df):
    # Calculate Simple Moving Average (SMA) over 20 periods (100 minutes)
    df['sma'] = df['Close'].rolling(window=20).mean()
    # Calculate Standard Deviation over 20 periods
    df['std_dev'] = df['Close'].rolling(window=20).std()
    # Define Bollinger Bands (2 standard deviations from SMA)
    df['lower_bb'] = df['sma'] - 2 * df['std_dev']
    df['upper_bb'] = df['sma'] + 2 * df['std_dev']
    # Initialize signals array with Hold (0) by default
    signals = np.zeros(len(df))
    # Assign Buy (1) when Close price falls below Lower Bollinger Band
    signals[(df['Close'] < df['lower_bb'])] = 1
    # Assign Sell (-1) when Close price rises above Upper Bollinger Band
    signals[(df['Close'] > df['upper_bb'])] = -1
    # Add signals to the original DataFrame
    df['signal'] = signals

'@
$ws.Range("C2").Value = $c2
$d2 = @'
content='**Analysis of Weaknesses and Improvements**
Based on the provided backtest results, the strategy's weaknesses are:
1. **Negative Total Return**: The strategy resulted in a loss of -0.83% over the 5-day period, indicating that the current buy/sell logic is not effective.
2. **High Max Drawdown**: A 0.60% max drawdown in a short period suggests that the strategy is prone to significant short-term losses.
3. **Low Sharpe Ratio**: A Sharpe Ratio of -0.71 implies that the strategy's returns are largely driven by risk rather than alpha.
4. **Frequent Trading**: 8 trades in 5 days may lead to high transaction costs, as evidenced by the $80 fee cost.
To address these weaknesses, the revised strategy will focus on:
1. **Improving risk management**: Reduce exposure to large losses.
2. **Enhancing return potential**: Identify more effective entry/exit points.
3. **Decreasing trading frequency**: Minimize transaction costs.
**Revised Strategy: "Mean Reversion with Volatility Guard"**
This strategy combines mean reversion with a volatility-based filter to reduce risk and improve returns.
**Python Code**

'@
$ws.Range("D2").Value = $d2
$e2 = @'

import numpy as np
import pandas as pd
def add_signal(df):
    """
    Adds a 'signal' column to the input DataFrame based on a mean reversion strategy with volatility guard.
    """
    # Calculate Short-Term Simple Moving Average (10 periods, 50 minutes)
    df['sma_short'] = df['Close'].rolling(window=10).mean()
    # Calculate Long-Term Simple Moving Average (30 periods, 150 minutes)
    df['sma_long'] = df['Close'].rolling(window=30).mean()
    # Calculate Standard Deviation over 20 periods (100 minutes)
    df['std_dev'] = df['Close'].rolling(window=20).std()
    # Volatility Guard: only trade when std_dev is below the 20-period average std_dev
    avg_std_dev = df['std_dev'].rolling(window=20).mean()
    volatility_guard = df['std_dev'] < avg_std_dev
    # Mean Reversion Logic
    # Buy when Close price falls below SMA (short) and SMA (short) is above SMA (long)
    buy_logic = (df['Close'] < df['sma_short']) & (df['sma_short'] > df['sma_long']) & volatility_guard
    # Sell when Close price rises above SMA (short) and SMA (short) is below SMA (long)
    sell_logic = (df['Close'] > df['sma_short']) & (df['sma_short'] < df['sma_long']) & volatility_guard
    # Initialize signals array with Hold (0) by default
    signals = np.zeros(len(df), dtype=int)
    # Assign Buy (1) and Sell (-1) signals
    signals[buy_logic] = 1
    signals[sell_logic] = -1
    # Add signals to the original DataFrame, ensuring same index
    df["signal"] = pd.Series(signals, index=df.index)

'@
$ws.Range("E2").Value = $e2
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").Value = -0.71
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = -29.21
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 70
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = 9

# ---- Row 3 (improved = Dual SMA crossover variant, new row) ----
$ws.Range("A3").Value = 45787.64977077832
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
$b3 = @'
content="### 1️⃣ Explanation: 
Mean Reversion Strategy for TSLA: Buy when 5-min Close price falls below its 20-period Simple Moving Average (SMA) by more than 1.5%, Sell when it rises above by more than 1.5%, otherwise Hold.

### 2️⃣ `add_signal` Function

'@
$ws.Range("B3").Value = $b3
$c3 = @'

def add_signal(df):
    # Calculate Simple Moving Average (SMA) over 20 periods
    df['sma_20'] = df['Close'].rolling(window=20).mean()

    # Calculate percentage difference between Close and SMA
    df['pct_diff'] = ((df['Close'] - df['sma_20']) / df['sma_20']) * 100

    # Initialize signal column with Hold (0) as default
    df['signal'] = 0

    # Set Buy (1) signal when Close is more than 1.5% below SMA
    df.loc[df['pct_diff'] < -1.5, 'signal'] = 1

    # Set Sell (-1) signal when Close is more than 1.5% above SMA
    df.loc[df['pct_diff'] > 1.5, 'signal'] = -1

'@
$ws.Range("C3").Value = $c3
$d3 = @'
content='**Analysis of Weaknesses and Improvements**

Based on the provided strategy and backtest results, the following weaknesses and potential improvements are identified:

1. **Overly Simple Logic**: The strategy relies solely on the percentage difference between the Close price and a 20-period SMA, which might not capture more nuanced market behaviors.
2. **High Trade Frequency**: With 8 trades over 5 trading days (approximately 1 trade per day), transaction costs (e.g., the $80 fee) might erode profits. Reducing trade frequency while maintaining or improving profitability is desirable.
3. **Limited Risk Management**: The strategy lacks explicit risk management techniques, such as position sizing based on volatility or maximum allowable loss per trade.
4. **Low Profitability**: The -0.06% total return over 5 days is underwhelming, suggesting the need for a more effective entry/exit logic.

**Revised Strategy: "Improved SMA Crossover with Volatility-Based Position Sizing"**

**Brief Explanation**:
This revised strategy aims to enhance profitability and reduce volatility by:

* Introducing a dual SMA crossover system for more informed entry/exit decisions.
* Incorporating a basic form of risk management through volatility-based position sizing (though, for simplicity, this aspect is simulated by adjusting the signal strength rather than actual position sizing, as the latter would require additional portfolio management logic not specified in the original task).

'@
$ws.Range("D3").Value = $d3
$e3 = @'

import pandas as pd
import numpy as np

def add_signal(df):
    """
    Revised strategy adding a 'signal' column to the DataFrame.
    Dual SMA Crossover with Simulated Volatility-Based Signal Strength
    """
    # Calculate Short and Long Simple Moving Averages
    df['sma_short'] = df['Close'].rolling(window=10).mean()
    df['sma_long'] = df['Close'].rolling(window=30).mean()
    # Calculate Volatility (Simple Measure: Close Price Standard Deviation over 20 periods)
    df['volatility'] = df['Close'].rolling(window=20).std()
    # Initialize signal column with Hold (0) as default
    df['signal'] = 0
    # Identify Buy and Sell Crossover Points
    buy_crossover = (df['sma_short'] > df['sma_long']) & (df['sma_short'].shift(1) <= df['sma_long'].shift(1))
    sell_crossover = (df['sma_short'] < df['sma_long']) & (df['sma_short'].shift(1) >= df['sma_long'].shift(1))
    # Assign Buy (1) and Sell (-1) Signals Based on Crossover, Adjusting Strength by Volatility
    df.loc[buy_crossover, 'signal'] = 1 / (1 + df['volatility'] / df['Close'])  # Simulated position sizing based on volatility
    df.loc[sell_crossover, 'signal'] = -1 / (1 + df['volatility'] / df['Close'])
    # Ensure signal Series uses the same index as df
    df['signal'] = pd.Series(df['signal'].values, index=df.index)

'@
$ws.Range("E3").Value = $e3
$ws.Range("J3").Value = 12.31
$ws.Range("K3").Value = 8
$ws.Range("L3").Value = 80
$ws.Range("M3").Value = 10
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = -22.27
$ws.Range("P3").Value = 5
$ws.Range("Q3").Value = 59.67
$ws.Range("R3").Value = 10
$ws.Range("S3").Value = 9
$ws.Range("T3").Value = 1.5
$ws.Range("U3").Value = 1.5

Write-Host "edit complete"
